$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity (D column) values to reflect stock reductions after sales
$ws.Range("D2").Value = 5
$ws.Range("D6").Value = 45
$ws.Range("D8").Value = 20
$ws.Range("D9").Value = 10
$ws.Range("D10").Value = 45

# Update creation_date (G11) to a new timestamp value (serial date number)
$ws.Range("G11").Value = 45818.70595449737
